# Update "想去人数" (interested-count) figures on the 展览 and 全部类型
# sheets to match the latest generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 28     # 抚州·BM次元盛典运动番only
$ws1.Range("F7").Value  = 1743   # 南昌·ACG CLUB动漫游戏嘉年华
$ws1.Range("F8").Value  = 42     # 南昌·原壤铁ONLY
$ws1.Range("F11").Value = 1893   # 南昌·CM02动漫游戏博览会
$ws1.Range("F12").Value = 137    # 信丰·端午节UPUP动漫展
$ws1.Range("F23").Value = 1022   # 江西·次元星河国风动漫游戏嘉年华
$ws1.Range("F25").Value = 319    # 南昌·幻梦境国际动漫游戏嘉年华1th
$ws1.Range("F28").Value = 275    # 赣州·第二届异次元动漫嘉年华

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 28     # 抚州·BM次元盛典运动番only
$ws4.Range("F7").Value  = 1743   # 南昌·ACG CLUB动漫游戏嘉年华
$ws4.Range("F9").Value  = 42     # 南昌·原壤铁ONLY
$ws4.Range("F12").Value = 1893   # 南昌·CM02动漫游戏博览会
$ws4.Range("F13").Value = 137    # 信丰·端午节UPUP动漫展
$ws4.Range("F24").Value = 1022   # 江西·次元星河国风动漫游戏嘉年华
$ws4.Range("F26").Value = 319    # 南昌·幻梦境国际动漫游戏嘉年华1th
$ws4.Range("F29").Value = 275    # 赣州·第二届异次元动漫嘉年华
